$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumericCell($ref, $value, $styleRef) {
    $cell = $ws.Range($ref)
    $cell.Value = $value
    $donor = $ws.Range($styleRef)
    $donor.Copy()
    $cell.PasteSpecial(-4122)
}

function Set-TextCell($ref, $text, $styleRef) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $donor = $ws.Range($styleRef)
    $donor.Copy()
    $cell.PasteSpecial(-4122)
}

# --- Header text updates (Volume/Number + report week date range) ---
$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Crime Complaints table updates (rows 16-28) ---
Set-TextCell "C16" "0" "A14"
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -57.142857142857
$ws.Range("J16").Value = 9
$ws.Range("K16").Value = -55.555555555555
$ws.Range("M16").Value = -69.230769230769
$ws.Range("N16").Value = -90
Set-NumericCell "C17" 2 "F16"
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 4
$ws.Range("H17").Value = -71.428571428571
$ws.Range("I17").Value = 11
$ws.Range("J17").Value = 23
$ws.Range("K17").Value = -52.173913043478
$ws.Range("L17").Value = -35.294117647058
$ws.Range("M17").Value = -21.428571428571
$ws.Range("N17").Value = -54.166666666666
Set-NumericCell "D18" 1 "F16"
Set-NumericCell "E18" -100 "L16"
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = -66.666666666666
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = -20
$ws.Range("M18").Value = -63.636363636363
$ws.Range("N18").Value = -95.180722891566
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -60
$ws.Range("F19").Value = 7
$ws.Range("G19").Value = 12
$ws.Range("H19").Value = -41.666666666666
$ws.Range("I19").Value = 15
$ws.Range("J19").Value = 26
$ws.Range("K19").Value = -42.307692307692
$ws.Range("L19").Value = -48.275862068965
$ws.Range("M19").Value = -16.666666666666
$ws.Range("N19").Value = -21.052631578947
$ws.Range("D20").Value = 1
Set-TextCell "F20" "0" "A14"
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -100
$ws.Range("J20").Value = 15
$ws.Range("K20").Value = -86.666666666666
$ws.Range("L20").Value = -60
$ws.Range("M20").Value = -33.333333333333
$ws.Range("N20").Value = -94.117647058823
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -71.428571428571
$ws.Range("F21").Value = 15
$ws.Range("G21").Value = 47
$ws.Range("H21").Value = -68.085106382978
$ws.Range("I21").Value = 37
$ws.Range("J21").Value = 78
$ws.Range("K21").Value = -52.564102564102
$ws.Range("L21").Value = -42.1875
$ws.Range("M21").Value = -38.333333333333
$ws.Range("N21").Value = -82.125603864734
Set-NumericCell "D22" 2 "F16"
Set-NumericCell "E22" -100 "L16"
Set-NumericCell "G22" 2 "F16"
Set-NumericCell "H22" -100 "L16"
Set-NumericCell "J22" 2 "F16"
Set-NumericCell "K22" -100 "L16"
Set-NumericCell "C23" 1 "F16"
Set-TextCell "D23" "0" "A14"
Set-TextCell "E23" "***.*" "A14"
Set-NumericCell "F23" 1 "F16"
$ws.Range("H23").Value = -66.666666666666
$ws.Range("I23").Value = 3
$ws.Range("K23").Value = -25
$ws.Range("L23").Value = -40
$ws.Range("M23").Value = -62.5
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -57.142857142857
$ws.Range("F24").Value = 20
$ws.Range("G24").Value = 36
$ws.Range("H24").Value = -44.444444444444
$ws.Range("I24").Value = 60
$ws.Range("J24").Value = 64
$ws.Range("K24").Value = -6.25
$ws.Range("L24").Value = -30.232558139534
$ws.Range("M24").Value = 33.333333333333
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -57.142857142857
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -52.941176470588
$ws.Range("I25").Value = 32
$ws.Range("J25").Value = 29
$ws.Range("K25").Value = 10.344827586206
$ws.Range("L25").Value = -34.693877551020
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 8
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = -12.5
$ws.Range("I26").Value = 32
$ws.Range("J26").Value = 36
$ws.Range("K26").Value = -11.111111111111
$ws.Range("L26").Value = 3.225806451612
$ws.Range("M26").Value = -27.272727272727
Set-NumericCell "C28" 1 "F16"
Set-NumericCell "F28" 1 "F16"
Set-TextCell "G28" "0" "A14"
Set-TextCell "H28" "***.*" "A14"
$ws.Range("I28").Value = 2
$ws.Range("K28").Value = -60
$ws.Range("L28").Value = -33.333333333333
